# Resolves #51: Add passive event check-ins and Mixtacular pull hat.
#
# The "Hats" sheet gains a new leading column ("Impl?") marking which hats
# have been implemented ("x") so far, shifting every existing column right
# by one. A handful of camo names on the "Camos" sheet are bolded, and the
# active selection on each sheet moves to reflect where the author was
# last working.

$wb = $excel.ActiveWorkbook

# --- Hats sheet: insert a new column A for the "Impl?" marker -------------
$wsHats = $wb.Worksheets.Item("Hats")
$wsHats.Columns("A:A").Insert()

$wsHats.Range("A1").Value2 = "Impl?"

# Hats that have been implemented so far: rows 2-18 (the "special friend"
# / donor hats) plus row 40 ("HAT_ALL_MIXERS").
$implementedRows = 2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,40
foreach ($r in $implementedRows) {
    $wsHats.Cells.Item($r, 1).Value2 = "x"
}

# --- Camos sheet: bold the camo-name cells for the newly tracked camos ----
$wsCamos = $wb.Worksheets.Item("Camos")
$boldRows = 5,7,13,15,16,17,22,23,24,25
foreach ($r in $boldRows) {
    $wsCamos.Cells.Item($r, 3).Font.Bold = $true
}

# --- Restore the author's last-used selection on each sheet ---------------
$wsEvents = $wb.Worksheets.Item("Event check-ins")

$wsHats.Activate()
$wsHats.Range("A41").Select()

$wsCamos.Activate()
$wsCamos.Range("G28").Select()

$wsEvents.Activate()
$wsEvents.Range("F17").Select()

$wsHats.Activate()

Write-Output "edit applied"
